# feat: add 2022-Q1 data
#
# The workbook currently has 4 sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计 (totals).
# 总计's last sheet holds the per-quarter summary (date / count / value).
#
# We need to:
#  1. Turn the existing "总计" sheet into the new "2022-Q1" sheet, holding the
#     per-fund holdings detail for the new quarter (same shape as the
#     2021-Q2 / 2021-Q3 / 2021-Q4 sheets).
#  2. Add a brand new "总计" sheet right after it, with the same summary
#     table as before plus a new first data row for 2022-Q1 (existing rows
#     shift down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet as "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"

# Remember a cell that already carries the shared header/index style (s=2)
# so we can stamp the same style onto the new header cells (E1:H1) without
# inventing a new style entry.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Wipe the old date/count/value rows (B2:D4) - headers (row1) and the A
# column index cells keep their formatting/values as-is.
$q1.Range("B2:D4").Clear()

# Headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund rows are stored as text (same pattern as the other quarter sheets),
# so force a text number format before writing the numeric-looking values.
$q1.Range("B2:G4").NumberFormat = "@"

$q1.Range("B2").Value = "006218"
$q1.Range("C2").Value = "富国生物医药科技混合A"
$q1.Range("D2").Value = "9.55"
$q1.Range("E2").Value = "87.25"
$q1.Range("F2").Value = "3.49"
$q1.Range("G2").Value = "0.3333"
$q1.Range("H2").Value = 10

$q1.Range("B3").Value = "100016"
$q1.Range("C3").Value = "富国天源沪港深平衡混合"
$q1.Range("D3").Value = "6.23"
$q1.Range("E3").Value = "72.29"
$q1.Range("F3").Value = "3.84"
$q1.Range("G3").Value = "0.2392"
$q1.Range("H3").Value = 4

$q1.Range("B4").Value = "011308"
$q1.Range("C4").Value = "富国生物医药科技混合C"
$q1.Range("D4").Value = "1.26"
$q1.Range("E4").Value = "87.25"
$q1.Range("F4").Value = "3.49"
$q1.Range("G4").Value = "0.0440"
$q1.Range("H4").Value = 10

# ---------------------------------------------------------------------
# Step 2: add the new "总计" sheet right after "2022-Q1"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Reuse the same shared style (s=2) for the header row and the A-column
# index cells, copying it from the styled cell we already have on $q1.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.62

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.3

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.58

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.01

# Restore the originally-active tab (2021-Q2) so adding sheets doesn't
# change which tab is selected.
$wb.Worksheets.Item(1).Activate()
